$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update the RF column (I) for rows 30 through 48 with the new value
$newRF = 9.012499999999999
for ($r = 30; $r -le 48; $r++) {
    $ws.Cells.Item($r, 9).Value = $newRF
}
